$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text describing the local file path to use a generic
# <SAPGUI-directory> placeholder instead of the hard-coded install path.
$ws.Range("B8").Value = "Image from a file (<SAPGUI-directory>\wwi\graphics\W_bio.bmp)"

# Refresh the (regenerated) picture shape names/ids.
$ws.Shapes.Item(1).Name = "1FCB857059FD1EDEBEAEE72CCDB39CB1"
$ws.Shapes.Item(2).Name = "1FCB857059FD1EDEBEAEE72CCDB3BCB1"
$ws.Shapes.Item(3).Name = "1FCB857059FD1EDEBEAEE72CCDB3DCB1"
